# Recalculate min atom distances: update residue-combination table rows 3-18
# and drop the now-unused row 19 (dimension shrinks from A1:D19 to A1:D18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Prefix with an apostrophe so Excel stores numeric-looking strings
    # (e.g. "780", "6016") as text, matching the original inlineStr cells.
    $ws.Range($addr).Value = "'" + $text
}

# Row 3
Set-TextCell "A3" "780"
Set-TextCell "C3" "780"

# Row 4
Set-TextCell "A4" "130, 1073"
$ws.Range("B4").Value = 2
Set-TextCell "C4" "130, 130"
Set-TextCell "D4" "5131, 4415"

# Row 5
Set-TextCell "A5" "780, 1073, 1105"
Set-TextCell "C5" "1105"
Set-TextCell "D5" "4994"

# Row 6
Set-TextCell "A6" "130, 455, 780"
$ws.Range("B6").Value = 2
Set-TextCell "C6" "130, 130"
Set-TextCell "D6" "5269, 6424"

# Row 7
Set-TextCell "A7" "455, 1073, 1105"
$ws.Range("B7").Value = 1
Set-TextCell "C7" "1105"
Set-TextCell "D7" "5399"

# Row 8
Set-TextCell "A8" "423, 748, 780, 1073"
$ws.Range("B8").Value = 2
Set-TextCell "C8" "780, 780"
Set-TextCell "D8" "5677, 5887"

# Row 9
Set-TextCell "A9" "423, 1073, 1105"
Set-TextCell "C9" "1105, 1105"
Set-TextCell "D9" "5331, 5433"

# Row 10
Set-TextCell "A10" "98, 130, 455, 748, 1073"
Set-TextCell "C10" "130"
Set-TextCell "D10" "6016"

# Row 11
Set-TextCell "A11" "98, 130, 748, 1073"
Set-TextCell "C11" "130"
Set-TextCell "D11" "5582"

# Row 12
Set-TextCell "A12" "130, 423, 748, 1073, SF"
Set-TextCell "D12" "6202"

# Row 13
Set-TextCell "A13" "98, 130, 780, 1073"
$ws.Range("B13").Value = 1
Set-TextCell "C13" "130"
Set-TextCell "D13" "6488"

# Row 14
Set-TextCell "A14" "130, 423, 780, 1073"
Set-TextCell "C14" "130"
Set-TextCell "D14" "6561"

# Row 15
Set-TextCell "A15" "130, 780, 780, 1073"
Set-TextCell "C15" "780"
Set-TextCell "D15" "6359"

# Row 16
Set-TextCell "A16" "98, 98, 130, 455, 780"
Set-TextCell "C16" "130"
Set-TextCell "D16" "6727"

# Row 17
Set-TextCell "A17" "98, 130, 423, 1073"
Set-TextCell "D17" "6670"

# Row 18
Set-TextCell "A18" "98, 98, 455, 780"
Set-TextCell "C18" "780"
Set-TextCell "D18" "6748"

# Row 19 no longer exists - clear it so the used range shrinks to A1:D18
$ws.Range("A19:D19").ClearContents()
